# Fixed issue #13: Permitir que en los ficheros de metadatos dos columnas se
# puedan relacionar para crear SKOS jerárquicos.
#
# Insert a new row right below the header row. This new row holds short
# "slug" identifiers (no iaest-measure:/sdmx-dimension: prefix) for each
# column; columns that should be related to another column (here M and N)
# simply reuse that other column's slug so the metadata loader can link
# them into a SKOS hierarchy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 2-4 down to make room for the new row 2.
$ws.Rows.Item(2).Insert()

$newRowValues = @(
    "sector-descripcion",
    "sector",
    "n-parados",
    "orden-nacionalidad",
    "provincia-nombre",
    "sexo",
    "ue28",
    "ue25-ue27-ue28",
    "ue27",
    "ue25",
    "sector-codigo",
    "codsect",
    "aragon",
    "provincia-codigo",
    "mes-y-ano"
)

for ($i = 0; $i -lt $newRowValues.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $newRowValues[$i]
}
